# Update "want to go" (想去人数) counts across all sheets to reflect
# the latest scrape output, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 345
$ws.Range("F5").Value  = 8101
$ws.Range("F9").Value  = 60
$ws.Range("F11").Value = 534
$ws.Range("F24").Value = 4575
$ws.Range("F26").Value = 48098
$ws.Range("F30").Value = 726
$ws.Range("F31").Value = 18
$ws.Range("F32").Value = 62
$ws.Range("F33").Value = 824
$ws.Range("F35").Value = 561
$ws.Range("F36").Value = 185
$ws.Range("F40").Value = 959
$ws.Range("F41").Value = 114
$ws.Range("F43").Value = 1044
$ws.Range("F45").Value = 86
$ws.Range("F46").Value = 72
$ws.Range("F47").Value = 18
$ws.Range("F48").Value = 2445

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 238
$ws.Range("F5").Value  = 50
$ws.Range("F6").Value  = 50
$ws.Range("F14").Value = 34
$ws.Range("F19").Value = 7272
$ws.Range("F30").Value = 10

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value  = 2194
$ws.Range("F5").Value  = 1471
$ws.Range("F9").Value  = 9231
$ws.Range("F10").Value = 1494
$ws.Range("F11").Value = 144
$ws.Range("F12").Value = 54

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 2194
$ws.Range("F4").Value  = 8101
$ws.Range("F5").Value  = 1471
$ws.Range("F7").Value  = 144
$ws.Range("F8").Value  = 54
$ws.Range("F10").Value = 60
$ws.Range("F11").Value = 534
$ws.Range("F13").Value = 238
$ws.Range("F14").Value = 50
$ws.Range("F18").Value = 4575
$ws.Range("F24").Value = 34
$ws.Range("F26").Value = 726
$ws.Range("F27").Value = 62
$ws.Range("F28").Value = 824
$ws.Range("F29").Value = 561
$ws.Range("F32").Value = 185
$ws.Range("F36").Value = 959
$ws.Range("F38").Value = 114
$ws.Range("F40").Value = 1044
$ws.Range("F43").Value = 86
$ws.Range("F45").Value = 72
$ws.Range("F46").Value = 18
$ws.Range("F48").Value = 2445
